$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7286480665206909
$ws.Range("B1").Value = 1.595165014266968
$ws.Range("C1").Value = 3.931602716445923
$ws.Range("D1").Value = 1.329577803611755
$ws.Range("E1").Value = 0.7497053146362305
